$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 12.87
$ws.Range("B3").Value = 6.367
$ws.Range("E5").Value = 13.254
$ws.Range("B14").Value = 6.645999999999999
$ws.Range("B21").Value = 6.179
$ws.Range("B23").Value = 6.665000000000001
$ws.Range("B25").Value = 6.159000000000001
